$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.459.50'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '3.492.56'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.71'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.65'
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("D7").Value = '3.491.12'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").Value = '4.092.51'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.119'
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").Value = '3.501.39'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '64.292.90'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.31'
$ws.Range("E18").Value = '  -8.17%  '
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.75'
$ws.Range("E21").Value = '  -3.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.93'
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.567'
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = '3.633.84'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.09'
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").Value = '  +5.54%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.54'
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").Value = '3.515.28'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.44'
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.29'
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.83'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.52'
$ws.Range("E40").Value = '  -4.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0781'
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.88'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("E47").Value = '  +1.84%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").Value = '2.476.08'
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.75'
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.904'
$ws.Range("E51").Value = '  +1.92%  '
